$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.319.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.867.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4820'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2803'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06502'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.864.88'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07465'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.47'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.076'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.275.48'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007586'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.105.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.293'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '220.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +14.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.157'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.302'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.964'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.456'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09327'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.303'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.021'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05041'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.204'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.91%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7469'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.72%  '
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01829'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9146'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.082'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.915'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4262'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.389'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1290'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '63.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.983'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.475'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.78'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05632'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.12%  '
